$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-21) were reshuffled (rows 5 and 6 stayed put).
# Set target values explicitly, row by row, columns A-F.

$data = @{
  2  = @(501, 9, 52, 30, 75, 45)
  3  = @(902, 1, 0, 0, 0, 0)
  4  = @(901, 16, 15, 45, 60, 60)
  5  = @(1001, 18, 30, 75, 60, 72)
  6  = @(301, 6, 45, 30, 60, 45)
  7  = @(401, 9, 48, 67, 75, 45)
  8  = @(701, 3, 90, 45, 97, 15)
  9  = @(601, 9, 60, 67, 60, 42)
  10 = @(801, 3, 67, 65, 52, 45)
  11 = @(101, 9, 30, 15, 60, 15)
  12 = @(1202, 2, 10, 10, 10, 10)
  13 = @(1201, 2, 10, 10, 10, 10)
  14 = @(201, 9, 30, 15, 45, 30)
  15 = @(1203, 3, 15, 15, 15, 15)
  16 = @(502, 0, 4, 0, 0, 0)
  17 = @(1, 0, 2, 2, 2, 2)
  18 = @(802, 0, 4, 5, 4, 0)
  19 = @(3, 0, 3, 3, 3, 3)
  20 = @(1101, 0, 15, 30, 30, 0)
  21 = @(2, 0, 2, 2, 2, 2)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $vals[$col - 1]
    }
}
